$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update header labels in row 1
$ws.Range("F1").Value = "First day - ENEM 2020"
$ws.Range("G1").Value = "Second day - ENEM 2020"

# Update F/G numeric values for rows 2-28 (rounded scores)
$values = @{
    2  = @(54.64, 57.31)
    3  = @(28.47, 32.23)
    4  = @(62.47, 63.49)
    5  = @(42.6, 47.21)
    6  = @(32.65, 35.79)
    7  = @(35.68, 38.91)
    8  = @(49.63, 51.91)
    9  = @(34.99, 37.91)
    10 = @(28.69, 30.61)
    11 = @(45.98, 50.04)
    12 = @(23.65, 25.89)
    13 = @(28.95, 30.91)
    14 = @(39.41, 41.86)
    15 = @(34.29, 36.74)
    16 = @(36.33, 38.51)
    17 = @(30.08, 32.42)
    18 = @(32.07, 34.34)
    19 = @(30.36, 34.05)
    20 = @(32.13, 35.15)
    21 = @(37.77, 41.03)
    22 = @(34.84, 37.49)
    23 = @(39.84, 42.94)
    24 = @(33.72, 36.64)
    25 = @(37.73, 41.14)
    26 = @(37.32, 40.32)
    27 = @(47.21, 49.85)
    28 = @(25.54, 27.86)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 6).Value = $pair[0]
    $ws.Cells.Item($row, 7).Value = $pair[1]
}
